# Auto-generated edit script: updates FFXIV leve-profit calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# market-board prices pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 936.92  # H2: 923.5 -> 936.92
$ws.Cells.Item(2, 9).Value = 844.381  # I2: 832.7273 -> 844.381
$ws.Cells.Item(2, 11).Value = 844.381  # K2: 832.7273 -> 844.381
$ws.Cells.Item(2, 13).Value = -731.381  # M2: -719.7273 -> -731.381
$ws.Cells.Item(28, 8).Value = 3716.5  # H28: 1900.04 -> 3716.5
$ws.Cells.Item(28, 9).Value = 3145.3635  # I28: 1562.5834 -> 3145.3635
$ws.Cells.Item(28, 11).Value = 3145.3635  # K28: 1562.5834 -> 3145.3635
$ws.Cells.Item(28, 13).Value = -2660.3635  # M28: -1077.5834 -> -2660.3635
$ws.Cells.Item(40, 8).Value = 6944.4443  # H40: 7937.5 -> 6944.4443
$ws.Cells.Item(40, 9).Value = 2000  # I40: 0 -> 2000
$ws.Cells.Item(40, 10).Value = 8357.143  # J40: 7937.5 -> 8357.143
$ws.Cells.Item(40, 11).Value = 2000  # K40: 0 -> 2000
$ws.Cells.Item(40, 12).Value = 8357.143  # L40: 7937.5 -> 8357.143
$ws.Cells.Item(40, 13).Value = -1825  # M40: None -> -1825
$ws.Cells.Item(40, 14).Value = -8707.143  # N40: -8287.5 -> -8707.143
$ws.Cells.Item(62, 8).Value = 2372.25  # H62: 2297.8 -> 2372.25
$ws.Cells.Item(62, 10).Value = 0  # J62: 2000 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 2000 -> 0
$ws.Cells.Item(62, 14).ClearContents()  # N62 was -3248
$ws.Cells.Item(65, 8).Value = 2372.25  # H65: 2297.8 -> 2372.25
$ws.Cells.Item(65, 10).Value = 0  # J65: 2000 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 10000 -> 0
$ws.Cells.Item(65, 14).ClearContents()  # N65 was -16240
$ws.Cells.Item(70, 8).Value = 2786.5  # H70: 2742.2856 -> 2786.5
$ws.Cells.Item(70, 9).Value = 2616  # I70: 2599.3333 -> 2616
$ws.Cells.Item(70, 10).Value = 2843.3333  # J70: 2849.5 -> 2843.3333
$ws.Cells.Item(70, 11).Value = 7848  # K70: 7797.999899999999 -> 7848
$ws.Cells.Item(70, 12).Value = 8529.999899999999  # L70: 8548.5 -> 8529.999899999999
$ws.Cells.Item(70, 13).Value = -7578  # M70: -7527.999899999999 -> -7578
$ws.Cells.Item(70, 14).Value = -9069.999899999999  # N70: -9088.5 -> -9069.999899999999
$ws.Cells.Item(73, 8).Value = 2786.5  # H73: 2742.2856 -> 2786.5
$ws.Cells.Item(73, 9).Value = 2616  # I73: 2599.3333 -> 2616
$ws.Cells.Item(73, 10).Value = 2843.3333  # J73: 2849.5 -> 2843.3333
$ws.Cells.Item(73, 11).Value = 7848  # K73: 7797.999899999999 -> 7848
$ws.Cells.Item(73, 12).Value = 8529.999899999999  # L73: 8548.5 -> 8529.999899999999
$ws.Cells.Item(73, 13).Value = -6912  # M73: -6861.999899999999 -> -6912
$ws.Cells.Item(73, 14).Value = -10401.9999  # N73: -10420.5 -> -10401.9999
$ws.Cells.Item(98, 8).Value = 11875  # H98: 11913.8 -> 11875
$ws.Cells.Item(98, 9).Value = 13343.75  # I98: 13392.25 -> 13343.75
$ws.Cells.Item(98, 11).Value = 13343.75  # K98: 13392.25 -> 13343.75
$ws.Cells.Item(98, 13).Value = -11845.75  # M98: -11894.25 -> -11845.75
$ws.Cells.Item(116, 8).Value = 5688.6665  # H116: 5695.9165 -> 5688.6665
$ws.Cells.Item(116, 9).Value = 5973.4287  # I116: 5837.625 -> 5973.4287
$ws.Cells.Item(116, 10).Value = 5290  # J116: 5412.5 -> 5290
$ws.Cells.Item(116, 11).Value = 5973.4287  # K116: 5837.625 -> 5973.4287
$ws.Cells.Item(116, 12).Value = 5290  # L116: 5412.5 -> 5290
$ws.Cells.Item(116, 13).Value = -2531.4287  # M116: -2395.625 -> -2531.4287
$ws.Cells.Item(116, 14).Value = -12174  # N116: -12296.5 -> -12174
$ws.Cells.Item(122, 8).Value = 11875  # H122: 11913.8 -> 11875
$ws.Cells.Item(122, 9).Value = 13343.75  # I122: 13392.25 -> 13343.75
$ws.Cells.Item(122, 11).Value = 40031.25  # K122: 40176.75 -> 40031.25
$ws.Cells.Item(122, 13).Value = -37581.25  # M122: -37726.75 -> -37581.25
$ws.Cells.Item(138, 8).Value = 1130720.8  # H138: 1073040 -> 1130720.8
$ws.Cells.Item(138, 9).Value = 2432.5557  # I138: 2536.75 -> 2432.5557
$ws.Cells.Item(138, 10).Value = 1493384.9  # J138: 1349298.9 -> 1493384.9
$ws.Cells.Item(138, 11).Value = 7297.6671  # K138: 7610.25 -> 7297.6671
$ws.Cells.Item(138, 12).Value = 4480154.699999999  # L138: 4047896.7 -> 4480154.699999999
$ws.Cells.Item(138, 13).Value = -2157.6671  # M138: -2470.25 -> -2157.6671
$ws.Cells.Item(138, 14).Value = -4490434.699999999  # N138: -4058176.7 -> -4490434.699999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2284.1765  # H2: 2165.3125 -> 2284.1765
$ws.Cells.Item(2, 9).Value = 2131.1428  # I2: 2109.7334 -> 2131.1428
$ws.Cells.Item(2, 10).Value = 2998.3333  # J2: 2999 -> 2998.3333
$ws.Cells.Item(2, 11).Value = 2131.1428  # K2: 2109.7334 -> 2131.1428
$ws.Cells.Item(2, 12).Value = 2998.3333  # L2: 2999 -> 2998.3333
$ws.Cells.Item(2, 13).Value = -2018.1428  # M2: -1996.7334 -> -2018.1428
$ws.Cells.Item(2, 14).Value = -3224.3333  # N2: -3225 -> -3224.3333
$ws.Cells.Item(32, 8).Value = 2456.682  # H32: 2482.5403 -> 2456.682
$ws.Cells.Item(32, 9).Value = 1725.8354  # I32: 1745.2693 -> 1725.8354
$ws.Cells.Item(32, 10).Value = 8871.888999999999  # J32: 8872.223 -> 8871.888999999999
$ws.Cells.Item(32, 11).Value = 1725.8354  # K32: 1745.2693 -> 1725.8354
$ws.Cells.Item(32, 12).Value = 8871.888999999999  # L32: 8872.223 -> 8871.888999999999
$ws.Cells.Item(32, 13).Value = -1438.8354  # M32: -1458.2693 -> -1438.8354
$ws.Cells.Item(32, 14).Value = -9445.888999999999  # N32: -9446.223 -> -9445.888999999999
$ws.Cells.Item(63, 8).Value = 4298  # H63: 4299.4116 -> 4298
$ws.Cells.Item(63, 10).Value = 6994.3335  # J63: 6998.3335 -> 6994.3335
$ws.Cells.Item(63, 12).Value = 6994.3335  # L63: 6998.3335 -> 6994.3335
$ws.Cells.Item(63, 14).Value = -8366.333500000001  # N63: -8370.333500000001 -> -8366.333500000001
$ws.Cells.Item(66, 8).Value = 4298  # H66: 4299.4116 -> 4298
$ws.Cells.Item(66, 10).Value = 6994.3335  # J66: 6998.3335 -> 6994.3335
$ws.Cells.Item(66, 12).Value = 34971.6675  # L66: 34991.6675 -> 34971.6675
$ws.Cells.Item(66, 14).Value = -41835.6675  # N66: -41855.6675 -> -41835.6675
$ws.Cells.Item(116, 8).Value = 2284.1765  # H116: 2165.3125 -> 2284.1765
$ws.Cells.Item(116, 9).Value = 2131.1428  # I116: 2109.7334 -> 2131.1428
$ws.Cells.Item(116, 10).Value = 2998.3333  # J116: 2999 -> 2998.3333
$ws.Cells.Item(116, 11).Value = 2131.1428  # K116: 2109.7334 -> 2131.1428
$ws.Cells.Item(116, 12).Value = 2998.3333  # L116: 2999 -> 2998.3333
$ws.Cells.Item(116, 13).Value = 162.8571999999999  # M116: 184.2665999999999 -> 162.8571999999999
$ws.Cells.Item(116, 14).Value = -7586.3333  # N116: -7587 -> -7586.3333
$ws.Cells.Item(132, 8).Value = 2843.0527  # H132: 2782.1 -> 2843.0527
$ws.Cells.Item(132, 9).Value = 2550.6  # I132: 2494.625 -> 2550.6
$ws.Cells.Item(132, 10).Value = 3939.75  # J132: 3932 -> 3939.75
$ws.Cells.Item(132, 11).Value = 7651.799999999999  # K132: 7483.875 -> 7651.799999999999
$ws.Cells.Item(132, 12).Value = 11819.25  # L132: 11796 -> 11819.25
$ws.Cells.Item(132, 13).Value = -5121.799999999999  # M132: -4953.875 -> -5121.799999999999
$ws.Cells.Item(132, 14).Value = -16879.25  # N132: -16856 -> -16879.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2284.1765  # H3: 2165.3125 -> 2284.1765
$ws.Cells.Item(3, 9).Value = 2131.1428  # I3: 2109.7334 -> 2131.1428
$ws.Cells.Item(3, 10).Value = 2998.3333  # J3: 2999 -> 2998.3333
$ws.Cells.Item(3, 11).Value = 2131.1428  # K3: 2109.7334 -> 2131.1428
$ws.Cells.Item(3, 12).Value = 2998.3333  # L3: 2999 -> 2998.3333
$ws.Cells.Item(3, 13).Value = -2017.1428  # M3: -1995.7334 -> -2017.1428
$ws.Cells.Item(3, 14).Value = -3226.3333  # N3: -3227 -> -3226.3333
$ws.Cells.Item(63, 8).Value = 0  # H63: 70271 -> 0
$ws.Cells.Item(63, 10).Value = 0  # J63: 70271 -> 0
$ws.Cells.Item(63, 12).Value = 0  # L63: 70271 -> 0
$ws.Cells.Item(63, 14).ClearContents()  # N63 was -71643
$ws.Cells.Item(66, 8).Value = 0  # H66: 70271 -> 0
$ws.Cells.Item(66, 10).Value = 0  # J66: 70271 -> 0
$ws.Cells.Item(66, 12).Value = 0  # L66: 210813 -> 0
$ws.Cells.Item(66, 14).ClearContents()  # N66 was -217677
$ws.Cells.Item(86, 8).Value = 6669.5  # H86: 2999.5 -> 6669.5
$ws.Cells.Item(86, 9).Value = 6669.5  # I86: 2999.5 -> 6669.5
$ws.Cells.Item(86, 11).Value = 6669.5  # K86: 2999.5 -> 6669.5
$ws.Cells.Item(86, 13).Value = -5546.5  # M86: -1876.5 -> -5546.5
$ws.Cells.Item(89, 8).Value = 6669.5  # H89: 2999.5 -> 6669.5
$ws.Cells.Item(89, 9).Value = 6669.5  # I89: 2999.5 -> 6669.5
$ws.Cells.Item(89, 11).Value = 33347.5  # K89: 14997.5 -> 33347.5
$ws.Cells.Item(89, 13).Value = -27731.5  # M89: -9381.5 -> -27731.5
$ws.Cells.Item(94, 8).Value = 6979.2  # H94: 6430.5 -> 6979.2
$ws.Cells.Item(94, 9).Value = 1224.25  # I94: 1299.6666 -> 1224.25
$ws.Cells.Item(94, 10).Value = 29999  # J94: 11561.333 -> 29999
$ws.Cells.Item(94, 11).Value = 1224.25  # K94: 1299.6666 -> 1224.25
$ws.Cells.Item(94, 12).Value = 29999  # L94: 11561.333 -> 29999
$ws.Cells.Item(94, 13).Value = -773.25  # M94: -848.6666 -> -773.25
$ws.Cells.Item(94, 14).Value = -30901  # N94: -12463.333 -> -30901

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3050.6667  # H31: 3135.3845 -> 3050.6667
$ws.Cells.Item(31, 9).Value = 3750  # I31: 5000 -> 3750
$ws.Cells.Item(31, 11).Value = 3750  # K31: 5000 -> 3750
$ws.Cells.Item(31, 13).Value = -3455  # M31: -4705 -> -3455
$ws.Cells.Item(34, 8).Value = 3050.6667  # H34: 3135.3845 -> 3050.6667
$ws.Cells.Item(34, 9).Value = 3750  # I34: 5000 -> 3750
$ws.Cells.Item(34, 11).Value = 3750  # K34: 5000 -> 3750
$ws.Cells.Item(34, 13).Value = -3548  # M34: -4798 -> -3548
$ws.Cells.Item(94, 8).Value = 2173.4546  # H94: 2156.4348 -> 2173.4546
$ws.Cells.Item(94, 9).Value = 2544.9285  # I94: 2580.6428 -> 2544.9285
$ws.Cells.Item(94, 10).Value = 1523.375  # J94: 1496.5555 -> 1523.375
$ws.Cells.Item(94, 11).Value = 2544.9285  # K94: 2580.6428 -> 2544.9285
$ws.Cells.Item(94, 12).Value = 1523.375  # L94: 1496.5555 -> 1523.375
$ws.Cells.Item(94, 13).Value = -2093.9285  # M94: -2129.6428 -> -2093.9285
$ws.Cells.Item(94, 14).Value = -2425.375  # N94: -2398.5555 -> -2425.375
$ws.Cells.Item(105, 8).Value = 2446.25  # H105: 2594.1428 -> 2446.25
$ws.Cells.Item(105, 10).Value = 4230.5  # J105: 5170.3335 -> 4230.5
$ws.Cells.Item(105, 12).Value = 4230.5  # L105: 5170.3335 -> 4230.5
$ws.Cells.Item(105, 14).Value = -7724.5  # N105: -8664.333500000001 -> -7724.5
$ws.Cells.Item(120, 8).Value = 21999  # H120: 0 -> 21999
$ws.Cells.Item(120, 10).Value = 21999  # J120: 0 -> 21999
$ws.Cells.Item(120, 12).Value = 21999  # L120: 0 -> 21999
$ws.Cells.Item(120, 14).Value = -29257  # N120: None -> -29257

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 620.25  # H14: 664.625 -> 620.25
$ws.Cells.Item(14, 9).Value = 620.25  # I14: 664.625 -> 620.25
$ws.Cells.Item(14, 11).Value = 1860.75  # K14: 1993.875 -> 1860.75
$ws.Cells.Item(14, 13).Value = -1687.75  # M14: -1820.875 -> -1687.75
$ws.Cells.Item(131, 8).Value = 20689.19  # H131: 15629.962 -> 20689.19
$ws.Cells.Item(131, 10).Value = 1641.1063  # J131: 1645.9706 -> 1641.1063
$ws.Cells.Item(131, 12).Value = 4923.3189  # L131: 4937.9118 -> 4923.3189
$ws.Cells.Item(131, 14).Value = -15003.3189  # N131: -15017.9118 -> -15003.3189

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9712.375  # H70: 10314.143 -> 9712.375
$ws.Cells.Item(73, 8).Value = 9712.375  # H73: 10314.143 -> 9712.375
$ws.Cells.Item(97, 8).Value = 749.9524  # H97: 762.5 -> 749.9524
$ws.Cells.Item(97, 9).Value = 750.2308  # I97: 771.1667 -> 750.2308
$ws.Cells.Item(97, 11).Value = 750.2308  # K97: 771.1667 -> 750.2308
$ws.Cells.Item(97, 13).Value = -254.2308  # M97: -275.1667 -> -254.2308
$ws.Cells.Item(122, 8).Value = 1497.7273  # H122: 1619.5555 -> 1497.7273
$ws.Cells.Item(122, 9).Value = 1610.8572  # I122: 1729.3334 -> 1610.8572
$ws.Cells.Item(122, 10).Value = 1299.75  # J122: 1400 -> 1299.75
$ws.Cells.Item(122, 11).Value = 4832.571599999999  # K122: 5188.0002 -> 4832.571599999999
$ws.Cells.Item(122, 12).Value = 3899.25  # L122: 4200 -> 3899.25
$ws.Cells.Item(122, 13).Value = -2382.571599999999  # M122: -2738.0002 -> -2382.571599999999
$ws.Cells.Item(122, 14).Value = -8799.25  # N122: -9100 -> -8799.25
$ws.Cells.Item(126, 8).Value = 5070.357  # H126: 5135.6 -> 5070.357
$ws.Cells.Item(126, 9).Value = 4030.111  # I126: 4058.875 -> 4030.111
$ws.Cells.Item(126, 10).Value = 6942.8  # J126: 6366.143 -> 6942.8
$ws.Cells.Item(126, 11).Value = 12090.333  # K126: 12176.625 -> 12090.333
$ws.Cells.Item(126, 12).Value = 20828.4  # L126: 19098.429 -> 20828.4
$ws.Cells.Item(126, 13).Value = -9620.332999999999  # M126: -9706.625 -> -9620.332999999999
$ws.Cells.Item(126, 14).Value = -25768.4  # N126: -24038.429 -> -25768.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2349.5652  # H46: 2197.6 -> 2349.5652
$ws.Cells.Item(46, 9).Value = 1535.4166  # I46: 1380.3572 -> 1535.4166
$ws.Cells.Item(46, 11).Value = 1535.4166  # K46: 1380.3572 -> 1535.4166
$ws.Cells.Item(46, 13).Value = -1347.4166  # M46: -1192.3572 -> -1347.4166
$ws.Cells.Item(61, 8).Value = 4243.357  # H61: 4815.1665 -> 4243.357
$ws.Cells.Item(61, 9).Value = 3367.3333  # I61: 3878.3 -> 3367.3333
$ws.Cells.Item(61, 11).Value = 3367.3333  # K61: 3878.3 -> 3367.3333
$ws.Cells.Item(61, 13).Value = -3165.3333  # M61: -3676.3 -> -3165.3333
$ws.Cells.Item(93, 8).Value = 1857.5714  # H93: 3750 -> 1857.5714
$ws.Cells.Item(93, 9).Value = 2040.8  # I93: 6000 -> 2040.8
$ws.Cells.Item(93, 10).Value = 1399.5  # J93: 1500 -> 1399.5
$ws.Cells.Item(93, 11).Value = 2040.8  # K93: 6000 -> 2040.8
$ws.Cells.Item(93, 12).Value = 1399.5  # L93: 1500 -> 1399.5
$ws.Cells.Item(93, 13).Value = -792.8  # M93: -4752 -> -792.8
$ws.Cells.Item(93, 14).Value = -3895.5  # N93: -3996 -> -3895.5
$ws.Cells.Item(113, 8).Value = 4243.357  # H113: 4815.1665 -> 4243.357
$ws.Cells.Item(113, 9).Value = 3367.3333  # I113: 3878.3 -> 3367.3333
$ws.Cells.Item(113, 11).Value = 3367.3333  # K113: 3878.3 -> 3367.3333
$ws.Cells.Item(113, 13).Value = -1197.3333  # M113: -1708.3 -> -1197.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 23288  # H42: 33772.4 -> 23288
$ws.Cells.Item(42, 9).Value = 26251  # I42: 41749.668 -> 26251
$ws.Cells.Item(42, 11).Value = 26251  # K42: 41749.668 -> 26251
$ws.Cells.Item(42, 13).Value = -25873  # M42: -41371.668 -> -25873
$ws.Cells.Item(107, 8).Value = 1374.875  # H107: 1500 -> 1374.875
$ws.Cells.Item(107, 9).Value = 999.5  # I107: 0 -> 999.5
$ws.Cells.Item(107, 11).Value = 2998.5  # K107: 0 -> 2998.5
$ws.Cells.Item(107, 13).Value = -1078.5  # M107: None -> -1078.5
$ws.Cells.Item(110, 8).Value = 149000  # H110: 142666.67 -> 149000
$ws.Cells.Item(110, 10).Value = 149000  # J110: 142666.67 -> 149000
$ws.Cells.Item(110, 12).Value = 149000  # L110: 142666.67 -> 149000
$ws.Cells.Item(110, 14).Value = -157180  # N110: -150846.67 -> -157180
$ws.Cells.Item(122, 8).Value = 5870.171  # H122: 5825.643 -> 5870.171
$ws.Cells.Item(122, 9).Value = 5807.4736  # I122: 5761.1284 -> 5807.4736
$ws.Cells.Item(122, 11).Value = 17422.4208  # K122: 17283.3852 -> 17422.4208
$ws.Cells.Item(122, 13).Value = -14972.4208  # M122: -14833.3852 -> -14972.4208

